$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 27, shifting existing rows 27-37 down to 28-38
$ws.Rows.Item(27).Insert()

# Fill the new row 27 with the same static columns as the other Jengibre rows,
# and the new data values from the diff.
$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 44777
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 100114007
$ws.Range("G27").Value = "Jengibre"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = 13000
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = 13600
$ws.Range("N27").Value = "$/caja 13 kilos"
$ws.Range("O27").Value = "Perú"
$ws.Range("P27").Value = 1046
$ws.Range("Q27").Value = 13
$ws.Range("R27").Value = "Hortaliza"
